# Monthly rollover update ("Actualizacion automatica"):
#  - "VENTA MENSUAL" sheet: the 4 rolling month columns (C:F) shift one
#    month to the left (oldest month drops off, a new empty month is
#    appended on the right) and the column headers follow the same shift.
#  - "VENTAS POR GRUPO" sheet: the per-category detail that belonged to the
#    month which just rolled out of view is cleared back to 0, and the
#    "<n> de 21" coverage counters in row 23 are reset to "0 de 21" for the
#    categories that lose their data.

$wb  = $excel.ActiveWorkbook
$wsGrupo   = $wb.Worksheets.Item("VENTAS POR GRUPO")
$wsMensual = $wb.Worksheets.Item("VENTA MENSUAL")

# ---------------------------------------------------------------------
# 1) "VENTAS POR GRUPO" - zero out the stale month's category figures
# ---------------------------------------------------------------------
$grupoCellsToZero = @(
    "D5", "M5", "O5",
    "L6",
    "M7",
    "M10",
    "D11", "H11", "M11",
    "G13", "H13", "M13", "O13",
    "M14",
    "H16", "I16", "M16",
    "L17", "M17",
    "D22", "H22", "I22", "M22"
)
foreach ($addr in $grupoCellsToZero) {
    $wsGrupo.Range($addr).Value = 0
}

# Row 23 coverage counters ("<n> de 21") -> reset to "0 de 21"
$grupoLabelCells = @("D23", "G23", "H23", "I23", "L23", "M23", "O23")
foreach ($addr in $grupoLabelCells) {
    $wsGrupo.Range($addr).Value = "0 de 21"
}

# ---------------------------------------------------------------------
# 2) "VENTA MENSUAL" - shift month columns C:F one position to the left
# ---------------------------------------------------------------------

# The column widths (authored in character units) ride along with the
# column-letter -> month remap: column C picks up column D's width, E
# picks up F's, and F (now the "new" last data column) reuses C's old
# width. Capture the "before" widths first since C gets overwritten
# before F needs it.
$widthC = $wsMensual.Columns.Item(3).ColumnWidth()
$widthD = $wsMensual.Columns.Item(4).ColumnWidth()
$widthF = $wsMensual.Columns.Item(6).ColumnWidth()

$wsMensual.Columns.Item(3).ColumnWidth = $widthD
$wsMensual.Columns.Item(5).ColumnWidth = $widthF
$wsMensual.Columns.Item(6).ColumnWidth = $widthC

# Header labels: junio/julio/agosto/septiembre -> julio/agosto/septiembre/octubre
$wsMensual.Range("C1").Value = "julio"
$wsMensual.Range("D1").Value = "agosto"
$wsMensual.Range("E1").Value = "septiembre"
$wsMensual.Range("F1").Value = "octubre"

# Data rows 2..23: new C = old D, new D = old E, new E = old F, new F = 0
# NOTE: ".Value" is an indexed/parameterised COM property, so it must be
# read by *calling* it - ".Value()" - not by bare member access (which
# would just hand back the property descriptor). Writing still uses plain
# assignment, ".Value = ...".
for ($row = 2; $row -le 23; $row++) {
    $oldD = $wsMensual.Cells.Item($row, 4).Value()
    $oldE = $wsMensual.Cells.Item($row, 5).Value()
    $oldF = $wsMensual.Cells.Item($row, 6).Value()

    $wsMensual.Cells.Item($row, 3).Value = $oldD
    $wsMensual.Cells.Item($row, 4).Value = $oldE
    $wsMensual.Cells.Item($row, 5).Value = $oldF
    $wsMensual.Cells.Item($row, 6).Value = 0
}
